$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.228.68'
$ws.Range("E2").Value = '  +2.20%  '
$ws.Range("D3").Value = '3.130.90'
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("D5").Value = "'576.97"
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").Value = "'180.62"
$ws.Range("E6").Value = '  +6.00%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.131.04'
$ws.Range("E8").Value = '  +1.99%  '
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("E10").Value = '  +2.45%  '
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("D12").Value = "'0.468"
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("E13").Value = '  +1.27%  '
$ws.Range("D14").Value = "'36.71"
$ws.Range("E14").Value = '  +2.85%  '
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = '68.120.25'
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("D17").Value = '3.653.05'
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("D19").Value = '3.132.15'
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").Value = "'16.62"
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("D21").Value = "'486.24"
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = "'7.80"
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").Value = "'0.696"
$ws.Range("E23").Value = '  +1.38%  '
$ws.Range("D24").Value = "'83.87"
$ws.Range("E24").Value = '  +1.27%  '
$ws.Range("D25").Value = "'12.96"
$ws.Range("E25").Value = '  +2.29%  '
$ws.Range("E26").Value = '  +5.82%  '
$ws.Range("D27").Value = "'10.57"
$ws.Range("E27").Value = '  +4.21%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = "'8.13"
$ws.Range("E29").Value = '  +4.21%  '
$ws.Range("E30").Value = '  +4.09%  '
$ws.Range("D31").Value = "'2.63"
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("D32").Value = "'28.19"
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("E33").Value = '  +0.97%  '
$ws.Range("D34").Value = '0.0₃0951'
$ws.Range("E34").Value = '  +3.91%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = "'48.79"
$ws.Range("E36").Value = '  +3.39%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").Value = "'5.65"
$ws.Range("E37").Value = '  +1.20%  '
$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").Value = "'0.955"
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("E39").Value = '  +8.17%  '
$ws.Range("E40").Value = '  +4.35%  '
$ws.Range("D41").Value = "'49.14"
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("E42").Value = '  +1.44%  '
$ws.Range("D43").Value = "'8.41"
$ws.Range("E43").Value = '  +1.09%  '
$ws.Range("E44").Value = '  +7.75%  '
$ws.Range("D45").Value = "'396.31"
$ws.Range("E45").Value = '  +8.17%  '
$ws.Range("D46").Value = '2.794.24'
$ws.Range("E46").Value = '  +1.41%  '
$ws.Range("D47").Value = "'27.08"
$ws.Range("E47").Value = '  +9.94%  '
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D49").Value = "'135.27"
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = "'2.35"
$ws.Range("E51").Value = '  +8.87%  '
